$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.946.22"
$ws.Range("E2").Value = "  -3.16%  "

$ws.Range("D3").Value = "2.917.68"
$ws.Range("E3").Value = "  -4.01%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'586.81"
$ws.Range("E5").Value = "  -1.42%  "

$ws.Range("D6").Value = "'145.02"
$ws.Range("E6").Value = "  -5.72%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -2.30%  "

$ws.Range("D9").Value = "2.914.02"
$ws.Range("E9").Value = "  -3.97%  "

$ws.Range("D10").Value = "'6.72"
$ws.Range("E10").Value = "  -1.55%  "

$ws.Range("E11").Value = "  -5.14%  "

$ws.Range("D12").Value = "'0.447"
$ws.Range("E12").Value = "  -3.73%  "

$ws.Range("E13").Value = "  -3.72%  "

$ws.Range("D14").Value = "'33.51"
$ws.Range("E14").Value = "  -6.51%  "

$ws.Range("E15").Value = "  +0.01%  "

$ws.Range("D16").Value = "3.402.24"
$ws.Range("E16").Value = "  -4.03%  "

$ws.Range("D17").Value = "60.919.51"
$ws.Range("E17").Value = "  -3.24%  "

$ws.Range("D18").Value = "'6.75"
$ws.Range("E18").Value = "  -4.84%  "

$ws.Range("D19").Value = "2.918.27"
$ws.Range("E19").Value = "  -4.05%  "

$ws.Range("D20").Value = "'428.31"
$ws.Range("E20").Value = "  -5.74%  "

$ws.Range("D21").Value = "'13.59"
$ws.Range("E21").Value = "  -5.10%  "

$ws.Range("D22").Value = "'0.681"
$ws.Range("E22").Value = "  -2.59%  "

$ws.Range("D23").Value = "'7.09"
$ws.Range("E23").Value = "  -5.84%  "

$ws.Range("D24").Value = "'80.68"
$ws.Range("E24").Value = "  -2.89%  "

$ws.Range("D25").Value = "'2.23"
$ws.Range("E25").Value = "  -3.19%  "

$ws.Range("D26").Value = "'10.70"
$ws.Range("E26").Value = "  -5.00%  "

$ws.Range("E27").Value = "  -3.39%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("E30").Value = "  -3.01%  "

$ws.Range("D31").Value = "'2.62"
$ws.Range("E31").Value = "  -3.04%  "

$ws.Range("D32").Value = "'2.16"
$ws.Range("E32").Value = "  -3.99%  "

$ws.Range("D33").Value = "'26.61"
$ws.Range("E33").Value = "  -3.72%  "

$ws.Range("E34").Value = "  -3.34%  "

$ws.Range("D35").Value = "0.0₃0872"
$ws.Range("E35").Value = "  +1.78%  "

$ws.Range("D36").Value = "'1.01"
$ws.Range("E36").Value = "  -3.29%  "

$ws.Range("D37").Value = "'5.64"
$ws.Range("E37").Value = "  -5.07%  "

$ws.Range("D38").Value = "'3.00"
$ws.Range("E38").Value = "  -5.28%  "

$ws.Range("D39").Value = "'0.126"
$ws.Range("E39").Value = "  -3.75%  "

$ws.Range("D40").Value = "'49.59"
$ws.Range("E40").Value = "  -1.61%  "

$ws.Range("E41").Value = "  -5.31%  "

$ws.Range("D42").Value = "'8.62"
$ws.Range("E42").Value = "  -5.84%  "

$ws.Range("D43").Value = "'0.297"
$ws.Range("E43").Value = "  -2.44%  "

$ws.Range("D44").Value = "'41.19"
$ws.Range("E44").Value = "  -5.55%  "

$ws.Range("D45").Value = "'380.28"
$ws.Range("E45").Value = "  -2.94%  "

$ws.Range("E46").Value = "  -3.11%  "

$ws.Range("D47").Value = "2.688.83"
$ws.Range("E47").Value = "  -1.26%  "

$ws.Range("D48").Value = "'132.85"
$ws.Range("E48").Value = "  +0.38%  "

$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("D50").Value = "'24.61"
$ws.Range("E50").Value = "  +0.13%  "

$ws.Range("D51").Value = "'0.106"
$ws.Range("E51").Value = "  -2.51%  "
